# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated K (column G) values per row, recalculated from Strike# to K
$kValues = @{
    2  = 2
    3  = 3
    4  = 4
    5  = 4
    6  = 2
    7  = 5
    8  = 2
    9  = 1
    10 = 2
    11 = 1
    12 = 1
    14 = 2
    16 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
